$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2969.1
$ws.Range("B3").Value = 2616.1
$ws.Range("B4").Value = 2605.4
$ws.Range("B5").Value = 2011.5
$ws.Range("B6").Value = 1897.7
$ws.Range("B7").Value = 1573.3
$ws.Range("B8").Value = 1515.8
$ws.Range("B9").Value = 1101.7
$ws.Range("B10").Value = 928.2

$ws.Range("E2").Select()
